# DNB Mastercard Demo -> generic "Sheet" with diversified merchant rows.
# Rebuilds rows 2-20 of column A (date), B (merchant text), E (Inn) and
# F (Ut) to match the target transaction list, widens the date number
# format to include a time component, and renames the worksheet tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Rename the worksheet tab -----------------------------------
$ws.Name = "Sheet"

# ---- 2. Extend the date-styled column (A) down to row 20 first, so
#        every row in the new range shares the same cell style as the
#        existing A2:A14 block before we touch the number format.
$ws.Range("A15:A20").NumberFormat = $ws.Range("A2").NumberFormat

# ---- 3. New transaction table (row, date-serial, merchant, inn, ut) -
$rows = @(
    @{ R = 2;  A = 45688; B = "SPOTIFY AB";              E = $null;  F = 129 },
    @{ R = 3;  A = 45687; B = "NETFLIX.COM";              E = $null;  F = 179 },
    @{ R = 4;  A = 45685; B = "KIWI MAJORSTUEN";           E = $null;  F = 456.8 },
    @{ R = 5;  A = 45684; B = "APOTEK 1 SINSEN";           E = $null;  F = 189 },
    @{ R = 6;  A = 45682; B = "REMA 1000 TORSHOV";         E = $null;  F = 892.3 },
    @{ R = 7;  A = 45681; B = "BURGER KING KARL JOHAN";    E = $null;  F = 159 },
    @{ R = 8;  A = 45679; B = "VINMONOPOLET OSLO S";       E = $null;  F = 675 },
    @{ R = 9;  A = 45677; B = "STARBUCKS KARL JOHAN";      E = $null;  F = 89 },
    @{ R = 10; A = 45675; B = "TANUM BOKHANDEL OSLO";      E = $null;  F = 349 },
    @{ R = 11; A = 45674; B = "GITHUB.COM";                E = $null;  F = 129 },
    @{ R = 12; A = 45672; B = "SAS EUROBONUS";             E = $null;  F = 2890 },
    @{ R = 13; A = 45671; B = "SATS GYM MAJORSTUEN";       E = $null;  F = 599 },
    @{ R = 14; A = 45669; B = "MENY BOGSTADVEIEN";         E = $null;  F = 567.45 },
    @{ R = 15; A = 45667; B = "POWER STORO";               E = $null;  F = 1299 },
    @{ R = 16; A = 45665; B = "MCDONALDS OSLO S";          E = $null;  F = 119 },
    @{ R = 17; A = 45664; B = "COOP EXTRA GRØNLAND";       E = $null;  F = 723.9 },
    @{ R = 18; A = 45662; B = "Innbetaling";               E = 15000;  F = $null },
    @{ R = 19; A = 45660; B = "XXL SPORT ALNA";            E = $null;  F = 1499 },
    @{ R = 20; A = 45659; B = "PEPPES PIZZA SOLLI";        E = $null;  F = 389 }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B

    if ($row.E -ne $null) {
        $ws.Range("E$r").Value = $row.E
    } else {
        $ws.Range("E$r").ClearContents()
    }

    if ($row.F -ne $null) {
        $ws.Range("F$r").Value = $row.F
    } else {
        $ws.Range("F$r").ClearContents()
    }
}

# ---- 4. Widen the custom date format to also carry a time component -
$ws.Range("A2:A20").NumberFormat = "yyyy-mm-dd h:mm:ss"
